$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 989.2963
$ws.Range("J17").Value = 989.2963
$ws.Range("L17").Value = 2967.8889
$ws.Range("N17").Value = -3303.8889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2784.1785
$ws.Range("I19").Value = 1374.5385
$ws.Range("J19").Value = 4005.8667
$ws.Range("K19").Value = 1374.5385
$ws.Range("L19").Value = 4005.8667
$ws.Range("M19").Value = -1199.5385
$ws.Range("N19").Value = -4355.8667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 705.6667
$ws.Range("I80").Value = 317.57144
$ws.Range("J80").Value = 865.4706
$ws.Range("K80").Value = 952.71432
$ws.Range("L80").Value = 2596.4118
$ws.Range("M80").Value = 45.28567999999996
$ws.Range("N80").Value = -4592.4118

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 705.6667
$ws.Range("I83").Value = 317.57144
$ws.Range("J83").Value = 865.4706
$ws.Range("K83").Value = 2858.14296
$ws.Range("L83").Value = 7789.2354
$ws.Range("M83").Value = 2133.85704
$ws.Range("N83").Value = -17773.2354

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2634.8333
$ws.Range("I94").Value = 2147.182
$ws.Range("K94").Value = 2147.182
$ws.Range("M94").Value = -1696.182

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1891.75
$ws.Range("I101").Value = 972.4286
$ws.Range("J101").Value = 3178.8
$ws.Range("K101").Value = 2917.2858
$ws.Range("L101").Value = 9536.400000000001
$ws.Range("M101").Value = -1295.2858
$ws.Range("N101").Value = -12780.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5499
$ws.Range("I106").Value = 5373.75
$ws.Range("K106").Value = 5373.75
$ws.Range("M106").Value = -4742.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4433.4287
$ws.Range("I138").Value = 3025.5715
$ws.Range("J138").Value = 5841.2856
$ws.Range("K138").Value = 9076.7145
$ws.Range("L138").Value = 17523.8568
$ws.Range("M138").Value = -3936.7145
$ws.Range("N138").Value = -27803.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1496.4
$ws.Range("I14").Value = 1370.5
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 1370.5
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = -1195.5
$ws.Range("N14").Value = -2350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3969.3594
$ws.Range("I32").Value = 3900.4644
$ws.Range("K32").Value = 3900.4644
$ws.Range("M32").Value = -3613.4644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3861.8914
$ws.Range("I61").Value = 3239.1333
$ws.Range("J61").Value = 5029.5625
$ws.Range("K61").Value = 3239.1333
$ws.Range("L61").Value = 5029.5625
$ws.Range("M61").Value = -3027.1333
$ws.Range("N61").Value = -5453.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2010.9048
$ws.Range("I74").Value = 1887.3636
$ws.Range("J74").Value = 2146.8
$ws.Range("K74").Value = 1887.3636
$ws.Range("L74").Value = 2146.8
$ws.Range("M74").Value = -1013.3636
$ws.Range("N74").Value = -3894.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2010.9048
$ws.Range("I77").Value = 1887.3636
$ws.Range("J77").Value = 2146.8
$ws.Range("K77").Value = 9436.817999999999
$ws.Range("L77").Value = 10734
$ws.Range("M77").Value = -5068.817999999999
$ws.Range("N77").Value = -19470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1486.6666
$ws.Range("I102").Value = 1468.7826
$ws.Range("K102").Value = 1468.7826
$ws.Range("M102").Value = 153.2174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1627.6296
$ws.Range("I122").Value = 1743.8695
$ws.Range("J122").Value = 959.25
$ws.Range("K122").Value = 5231.6085
$ws.Range("L122").Value = 2877.75
$ws.Range("M122").Value = -2781.6085
$ws.Range("N122").Value = -7777.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3026.9722
$ws.Range("I132").Value = 2941.1292
$ws.Range("J132").Value = 3559.2
$ws.Range("K132").Value = 8823.3876
$ws.Range("L132").Value = 10677.6
$ws.Range("M132").Value = -6293.3876
$ws.Range("N132").Value = -15737.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3861.8914
$ws.Range("I136").Value = 3239.1333
$ws.Range("J136").Value = 5029.5625
$ws.Range("K136").Value = 9717.3999
$ws.Range("L136").Value = 15088.6875
$ws.Range("M136").Value = -7167.3999
$ws.Range("N136").Value = -20188.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 80957
$ws.Range("J139").Value = 80957
$ws.Range("L139").Value = 80957
$ws.Range("N139").Value = -91237

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5448.3774
$ws.Range("I134").Value = 4665.773
$ws.Range("K134").Value = 13997.319
$ws.Range("M134").Value = -11462.319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9516.714
$ws.Range("I31").Value = 4183.2
$ws.Range("J31").Value = 11183.4375
$ws.Range("K31").Value = 4183.2
$ws.Range("L31").Value = 11183.4375
$ws.Range("M31").Value = -3888.2
$ws.Range("N31").Value = -11773.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9516.714
$ws.Range("I34").Value = 4183.2
$ws.Range("J34").Value = 11183.4375
$ws.Range("K34").Value = 4183.2
$ws.Range("L34").Value = 11183.4375
$ws.Range("M34").Value = -3981.2
$ws.Range("N34").Value = -11587.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1600.2778
$ws.Range("I94").Value = 1718.6364
$ws.Range("K94").Value = 1718.6364
$ws.Range("M94").Value = -1267.6364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1303.6666
$ws.Range("I107").Value = 949.7273
$ws.Range("J107").Value = 2277
$ws.Range("K107").Value = 949.7273
$ws.Range("L107").Value = 2277
$ws.Range("M107").Value = 970.2727
$ws.Range("N107").Value = -6117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4112.7407
$ws.Range("I122").Value = 3933.85
$ws.Range("J122").Value = 4623.857
$ws.Range("K122").Value = 11801.55
$ws.Range("L122").Value = 13871.571
$ws.Range("M122").Value = -9351.549999999999
$ws.Range("N122").Value = -18771.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2670739.8
$ws.Range("I11").Value = 3129609.5
$ws.Range("J11").Value = 1753000
$ws.Range("K11").Value = 9388828.5
$ws.Range("L11").Value = 5259000
$ws.Range("M11").Value = -9388688.5
$ws.Range("N11").Value = -5259280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 4456.3335
$ws.Range("I16").Value = 2099.6667
$ws.Range("J16").Value = 5634.6665
$ws.Range("K16").Value = 6299.000100000001
$ws.Range("L16").Value = 16903.9995
$ws.Range("M16").Value = -6126.000100000001
$ws.Range("N16").Value = -17249.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3164.4285
$ws.Range("I68").Value = 2625.75
$ws.Range("J68").Value = 3291.1765
$ws.Range("K68").Value = 7877.25
$ws.Range("L68").Value = 9873.529500000001
$ws.Range("M68").Value = -7066.25
$ws.Range("N68").Value = -11495.5295

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3164.4285
$ws.Range("I71").Value = 2625.75
$ws.Range("J71").Value = 3291.1765
$ws.Range("K71").Value = 23631.75
$ws.Range("L71").Value = 29620.5885
$ws.Range("M71").Value = -19575.75
$ws.Range("N71").Value = -37732.5885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4087.8965
$ws.Range("I132").Value = 3309.52
$ws.Range("J132").Value = 8952.75
$ws.Range("K132").Value = 9928.559999999999
$ws.Range("L132").Value = 26858.25
$ws.Range("M132").Value = -7398.559999999999
$ws.Range("N132").Value = -31918.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3793014.8
$ws.Range("I40").Value = 4633840.5
$ws.Range("J40").Value = 9299
$ws.Range("K40").Value = 4633840.5
$ws.Range("L40").Value = 9299
$ws.Range("M40").Value = -4633704.5
$ws.Range("N40").Value = -9571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 318
$ws.Range("I100").Value = 318
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 636
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -95

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1332.4706
$ws.Range("I122").Value = 1276.6666
$ws.Range("K122").Value = 3829.9998
$ws.Range("M122").Value = -1379.9998
